$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 8.467421333333334
$ws.Cells.Item(2, 8).Value = 25.402264
$ws.Cells.Item(2, 9).Value = 0.2732469334691616
$ws.Cells.Item(2, 10).Value = 0.312800300005396
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 37.89292233333333
$ws.Cells.Item(2, 14).Value = 113.678767
$ws.Cells.Item(2, 15).Value = 0.1675903872431219
$ws.Cells.Item(2, 16).Value = 0.1770445447021447
$ws.Cells.Item(2, 17).Value = 320.8553389476098
$ws.Cells.Item(2, 18).Value = 2887.698050528488
$ws.Cells.Item(2, 19).Value = 0.04579355939309237
$ws.Cells.Item(2, 20).Value = 0.0553795866971496
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 8.467421333333334
$ws.Cells.Item(3, 8).Value = 25.402264
$ws.Cells.Item(3, 9).Value = 0.2732469334691616
$ws.Cells.Item(3, 10).Value = 0.312800300005396
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 49.58946566666666
$ws.Cells.Item(3, 14).Value = 148.768397
$ws.Cells.Item(3, 15).Value = 0.2193211091282201
$ws.Cells.Item(3, 16).Value = 0.2316935150513456
$ws.Cells.Item(3, 17).Value = 419.8948994945342
$ws.Cells.Item(3, 18).Value = 3779.054095450808
$ws.Cells.Item(3, 19).Value = 0.05992882051434147
$ws.Cells.Item(3, 20).Value = 0.07247380101736564
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 8.467421333333334
$ws.Cells.Item(4, 8).Value = 25.402264
$ws.Cells.Item(4, 9).Value = 0.2732469334691616
$ws.Cells.Item(4, 10).Value = 0.312800300005396
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 58.78086233333334
$ws.Cells.Item(4, 14).Value = 176.342587
$ws.Cells.Item(4, 15).Value = 0.2599722289632498
$ws.Cells.Item(4, 16).Value = 0.2746378576309976
$ws.Cells.Item(4, 17).Value = 497.7223277129966
$ws.Cells.Item(4, 18).Value = 4479.500949416969
$ws.Cells.Item(4, 19).Value = 0.07103661435135075
$ws.Cells.Item(4, 20).Value = 0.08590680425981527
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 8.467421333333334
$ws.Cells.Item(5, 8).Value = 25.402264
$ws.Cells.Item(5, 9).Value = 0.2732469334691616
$ws.Cells.Item(5, 10).Value = 0.312800300005396
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 43.61929799999999
$ws.Cells.Item(5, 14).Value = 130.857894
$ws.Cells.Item(5, 15).Value = 0.1929166343727092
$ws.Cells.Item(5, 16).Value = 0.2037995034192402
$ws.Cells.Item(5, 17).Value = 369.342974430224
$ws.Cells.Item(5, 18).Value = 3324.086769872016
$ws.Cells.Item(5, 19).Value = 0.05271387875753425
$ws.Cells.Item(5, 20).Value = 0.06374854581048907
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 8.467421333333334
$ws.Cells.Item(6, 8).Value = 25.402264
$ws.Cells.Item(6, 9).Value = 0.2732469334691616
$ws.Cells.Item(6, 10).Value = 0.312800300005396
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 13).Value = 36.221842
$ws.Cells.Item(6, 14).Value = 72.443684
$ws.Cells.Item(6, 15).Value = 0.160199640292699
$ws.Cells.Item(6, 16).Value = 0.112824579196272
$ws.Cells.Item(6, 17).Value = 306.7055976834294
$ws.Cells.Item(6, 18).Value = 1840.233586100576
$ws.Cells.Item(6, 19).Value = 0.04377406045284273
$ws.Cells.Item(6, 20).Value = 0.03529156222057642
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 0.3343473333333333
$ws.Cells.Item(7, 8).Value = 1.003042
$ws.Cells.Item(7, 9).Value = 0.01078951666043526
$ws.Cells.Item(7, 10).Value = 0.01235133366529898
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 37.89292233333333
$ws.Cells.Item(7, 14).Value = 113.678767
$ws.Cells.Item(7, 15).Value = 0.1675903872431219
$ws.Cells.Item(7, 16).Value = 0.1770445447021447
$ws.Cells.Item(7, 17).Value = 12.66939753435711
$ws.Cells.Item(7, 18).Value = 114.024577809214
$ws.Cells.Item(7, 19).Value = 0.00180821927528846
$ws.Cells.Item(7, 20).Value = 0.00218673624523713
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 0.3343473333333333
$ws.Cells.Item(8, 8).Value = 1.003042
$ws.Cells.Item(8, 9).Value = 0.01078951666043526
$ws.Cells.Item(8, 10).Value = 0.01235133366529898
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 49.58946566666666
$ws.Cells.Item(8, 14).Value = 148.768397
$ws.Cells.Item(8, 15).Value = 0.2193211091282201
$ws.Cells.Item(8, 16).Value = 0.2316935150513456
$ws.Cells.Item(8, 17).Value = 16.58010560707489
$ws.Cells.Item(8, 18).Value = 149.220950463674
$ws.Cells.Item(8, 19).Value = 0.002366368760924069
$ws.Cells.Item(8, 20).Value = 0.002861723912485141
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 0.3343473333333333
$ws.Cells.Item(9, 8).Value = 1.003042
$ws.Cells.Item(9, 9).Value = 0.01078951666043526
$ws.Cells.Item(9, 10).Value = 0.01235133366529898
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 58.78086233333334
$ws.Cells.Item(9, 14).Value = 176.342587
$ws.Cells.Item(9, 15).Value = 0.2599722289632498
$ws.Cells.Item(9, 16).Value = 0.2746378576309976
$ws.Cells.Item(9, 17).Value = 19.65322457218378
$ws.Cells.Item(9, 18).Value = 176.879021149654
$ws.Cells.Item(9, 19).Value = 0.002804974695649472
$ws.Cells.Item(9, 20).Value = 0.003392143816723329
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 0.3343473333333333
$ws.Cells.Item(10, 8).Value = 1.003042
$ws.Cells.Item(10, 9).Value = 0.01078951666043526
$ws.Cells.Item(10, 10).Value = 0.01235133366529898
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 43.61929799999999
$ws.Cells.Item(10, 14).Value = 130.857894
$ws.Cells.Item(10, 15).Value = 0.1929166343727092
$ws.Cells.Item(10, 16).Value = 0.2037995034192402
$ws.Cells.Item(10, 17).Value = 14.583995968172
$ws.Cells.Item(10, 18).Value = 131.255963713548
$ws.Cells.Item(10, 19).Value = 0.002081477240639443
$ws.Cells.Item(10, 20).Value = 0.002517195667553277
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 0.3343473333333333
$ws.Cells.Item(11, 8).Value = 1.003042
$ws.Cells.Item(11, 9).Value = 0.01078951666043526
$ws.Cells.Item(11, 10).Value = 0.01235133366529898
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 13).Value = 36.221842
$ws.Cells.Item(11, 14).Value = 72.443684
$ws.Cells.Item(11, 15).Value = 0.160199640292699
$ws.Cells.Item(11, 16).Value = 0.112824579196272
$ws.Cells.Item(11, 17).Value = 12.11067628112133
$ws.Cells.Item(11, 18).Value = 72.664057686728
$ws.Cells.Item(11, 19).Value = 0.001728476687933811
$ws.Cells.Item(11, 20).Value = 0.001393534023300105
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 4.140032333333333
$ws.Cells.Item(12, 8).Value = 12.420097
$ws.Cells.Item(12, 9).Value = 0.1336004309946363
$ws.Cells.Item(12, 10).Value = 0.1529395201819853
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 37.89292233333333
$ws.Cells.Item(12, 14).Value = 113.678767
$ws.Cells.Item(12, 15).Value = 0.1675903872431219
$ws.Cells.Item(12, 16).Value = 0.1770445447021447
$ws.Cells.Item(12, 17).Value = 156.8779236644887
$ws.Cells.Item(12, 18).Value = 1411.901312980399
$ws.Cells.Item(12, 19).Value = 0.02239014796623908
$ws.Cells.Item(12, 20).Value = 0.02707710771758405
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 4.140032333333333
$ws.Cells.Item(13, 8).Value = 12.420097
$ws.Cells.Item(13, 9).Value = 0.1336004309946363
$ws.Cells.Item(13, 10).Value = 0.1529395201819853
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 49.58946566666666
$ws.Cells.Item(13, 14).Value = 148.768397
$ws.Cells.Item(13, 15).Value = 0.2193211091282201
$ws.Cells.Item(13, 16).Value = 0.2316935150513456
$ws.Cells.Item(13, 17).Value = 205.3019912527232
$ws.Cells.Item(13, 18).Value = 1847.717921274509
$ws.Cells.Item(13, 19).Value = 0.02930139470575185
$ws.Cells.Item(13, 20).Value = 0.03543509502123039
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 4.140032333333333
$ws.Cells.Item(14, 8).Value = 12.420097
$ws.Cells.Item(14, 9).Value = 0.1336004309946363
$ws.Cells.Item(14, 10).Value = 0.1529395201819853
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 58.78086233333334
$ws.Cells.Item(14, 14).Value = 176.342587
$ws.Cells.Item(14, 15).Value = 0.2599722289632498
$ws.Cells.Item(14, 16).Value = 0.2746378576309976
$ws.Cells.Item(14, 17).Value = 243.3546706412154
$ws.Cells.Item(14, 18).Value = 2190.192035770939
$ws.Cells.Item(14, 19).Value = 0.03473240183612643
$ws.Cells.Item(14, 20).Value = 0.04200298216989316
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 4.140032333333333
$ws.Cells.Item(15, 8).Value = 12.420097
$ws.Cells.Item(15, 9).Value = 0.1336004309946363
$ws.Cells.Item(15, 10).Value = 0.1529395201819853
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 43.61929799999999
$ws.Cells.Item(15, 14).Value = 130.857894
$ws.Cells.Item(15, 15).Value = 0.1929166343727092
$ws.Cells.Item(15, 16).Value = 0.2037995034192402
$ws.Cells.Item(15, 17).Value = 180.585304077302
$ws.Cells.Item(15, 18).Value = 1625.267736695718
$ws.Cells.Item(15, 19).Value = 0.02577374549822861
$ws.Cells.Item(15, 20).Value = 0.03116899826626547
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 4.140032333333333
$ws.Cells.Item(16, 8).Value = 12.420097
$ws.Cells.Item(16, 9).Value = 0.1336004309946363
$ws.Cells.Item(16, 10).Value = 0.1529395201819853
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 13).Value = 36.221842
$ws.Cells.Item(16, 14).Value = 72.443684
$ws.Cells.Item(16, 15).Value = 0.160199640292699
$ws.Cells.Item(16, 16).Value = 0.112824579196272
$ws.Cells.Item(16, 17).Value = 149.9595970528913
$ws.Cells.Item(16, 18).Value = 899.7575823173481
$ws.Cells.Item(16, 19).Value = 0.02140274098829028
$ws.Cells.Item(16, 20).Value = 0.01725533700701223
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 6.291073
$ws.Cells.Item(17, 8).Value = 18.873219
$ws.Cells.Item(17, 9).Value = 0.2030153381778063
$ws.Cells.Item(17, 10).Value = 0.2324024569332694
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 37.89292233333333
$ws.Cells.Item(17, 14).Value = 113.678767
$ws.Cells.Item(17, 15).Value = 0.1675903872431219
$ws.Cells.Item(17, 16).Value = 0.1770445447021447
$ws.Cells.Item(17, 17).Value = 238.3871405823303
$ws.Cells.Item(17, 18).Value = 2145.484265240973
$ws.Cells.Item(17, 19).Value = 0.03402341914151192
$ws.Cells.Item(17, 20).Value = 0.04114558717541047
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 6.291073
$ws.Cells.Item(18, 8).Value = 18.873219
$ws.Cells.Item(18, 9).Value = 0.2030153381778063
$ws.Cells.Item(18, 10).Value = 0.2324024569332694
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 49.58946566666666
$ws.Cells.Item(18, 14).Value = 148.768397
$ws.Cells.Item(18, 15).Value = 0.2193211091282201
$ws.Cells.Item(18, 16).Value = 0.2316935150513456
$ws.Cells.Item(18, 17).Value = 311.9709485399936
$ws.Cells.Item(18, 18).Value = 2807.738536859943
$ws.Cells.Item(18, 19).Value = 0.04452554913919716
$ws.Cells.Item(18, 20).Value = 0.05384614215343814
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 6.291073
$ws.Cells.Item(19, 8).Value = 18.873219
$ws.Cells.Item(19, 9).Value = 0.2030153381778063
$ws.Cells.Item(19, 10).Value = 0.2324024569332694
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 13).Value = 58.78086233333334
$ws.Cells.Item(19, 14).Value = 176.342587
$ws.Cells.Item(19, 15).Value = 0.2599722289632498
$ws.Cells.Item(19, 16).Value = 0.2746378576309976
$ws.Cells.Item(19, 17).Value = 369.7946959419504
$ws.Cells.Item(19, 18).Value = 3328.152263477553
$ws.Cells.Item(19, 19).Value = 0.05277834997981225
$ws.Cells.Item(19, 20).Value = 0.06382651288033328
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 7).Value = 6.291073
$ws.Cells.Item(20, 8).Value = 18.873219
$ws.Cells.Item(20, 9).Value = 0.2030153381778063
$ws.Cells.Item(20, 10).Value = 0.2324024569332694
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 13).Value = 43.61929799999999
$ws.Cells.Item(20, 14).Value = 130.857894
$ws.Cells.Item(20, 15).Value = 0.1929166343727092
$ws.Cells.Item(20, 16).Value = 0.2037995034192402
$ws.Cells.Item(20, 17).Value = 274.412187926754
$ws.Cells.Item(20, 18).Value = 2469.709691340785
$ws.Cells.Item(20, 19).Value = 0.03916503576729978
$ws.Cells.Item(20, 20).Value = 0.04736350531641165
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 7).Value = 6.291073
$ws.Cells.Item(21, 8).Value = 18.873219
$ws.Cells.Item(21, 9).Value = 0.2030153381778063
$ws.Cells.Item(21, 10).Value = 0.2324024569332694
$ws.Cells.Item(21, 11).Value = 2
$ws.Cells.Item(21, 13).Value = 36.221842
$ws.Cells.Item(21, 14).Value = 72.443684
$ws.Cells.Item(21, 15).Value = 0.160199640292699
$ws.Cells.Item(21, 16).Value = 0.112824579196272
$ws.Cells.Item(21, 17).Value = 227.874252216466
$ws.Cells.Item(21, 18).Value = 1367.245513298796
$ws.Cells.Item(21, 19).Value = 0.03252298414998521
$ws.Cells.Item(21, 20).Value = 0.02622070940767583
$ws.Cells.Item(22, 5).Value = 2
$ws.Cells.Item(22, 7).Value = 11.755292
$ws.Cells.Item(22, 8).Value = 23.510584
$ws.Cells.Item(22, 9).Value = 0.3793477806979606
$ws.Cells.Item(22, 10).Value = 0.2895063892140504
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 13).Value = 37.89292233333333
$ws.Cells.Item(22, 14).Value = 113.678767
$ws.Cells.Item(22, 15).Value = 0.1675903872431219
$ws.Cells.Item(22, 16).Value = 0.1770445447021447
$ws.Cells.Item(22, 17).Value = 445.4423667616546
$ws.Cells.Item(22, 18).Value = 2672.654200569928
$ws.Cells.Item(22, 19).Value = 0.06357504146699011
$ws.Cells.Item(22, 20).Value = 0.05125552686676346
$ws.Cells.Item(23, 5).Value = 2
$ws.Cells.Item(23, 7).Value = 11.755292
$ws.Cells.Item(23, 8).Value = 23.510584
$ws.Cells.Item(23, 9).Value = 0.3793477806979606
$ws.Cells.Item(23, 10).Value = 0.2895063892140504
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 13).Value = 49.58946566666666
$ws.Cells.Item(23, 14).Value = 148.768397
$ws.Cells.Item(23, 15).Value = 0.2193211091282201
$ws.Cells.Item(23, 16).Value = 0.2316935150513456
$ws.Cells.Item(23, 17).Value = 582.9386490356412
$ws.Cells.Item(23, 18).Value = 3497.631894213847
$ws.Cells.Item(23, 19).Value = 0.08319897600800551
$ws.Cells.Item(23, 20).Value = 0.06707675294682632
$ws.Cells.Item(24, 5).Value = 2
$ws.Cells.Item(24, 7).Value = 11.755292
$ws.Cells.Item(24, 8).Value = 23.510584
$ws.Cells.Item(24, 9).Value = 0.3793477806979606
$ws.Cells.Item(24, 10).Value = 0.2895063892140504
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 13).Value = 58.78086233333334
$ws.Cells.Item(24, 14).Value = 176.342587
$ws.Cells.Item(24, 15).Value = 0.2599722289632498
$ws.Cells.Item(24, 16).Value = 0.2746378576309976
$ws.Cells.Item(24, 17).Value = 690.9862007401347
$ws.Cells.Item(24, 18).Value = 4145.917204440808
$ws.Cells.Item(24, 19).Value = 0.09861988810031089
$ws.Cells.Item(24, 20).Value = 0.07950941450423255
$ws.Cells.Item(25, 5).Value = 2
$ws.Cells.Item(25, 7).Value = 11.755292
$ws.Cells.Item(25, 8).Value = 23.510584
$ws.Cells.Item(25, 9).Value = 0.3793477806979606
$ws.Cells.Item(25, 10).Value = 0.2895063892140504
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(25, 13).Value = 43.61929799999999
$ws.Cells.Item(25, 14).Value = 130.857894
$ws.Cells.Item(25, 15).Value = 0.1929166343727092
$ws.Cells.Item(25, 16).Value = 0.2037995034192402
$ws.Cells.Item(25, 17).Value = 512.7575848250159
$ws.Cells.Item(25, 18).Value = 3076.545508950096
$ws.Cells.Item(25, 19).Value = 0.07318249710900715
$ws.Cells.Item(25, 20).Value = 0.05900125835852076
$ws.Cells.Item(26, 5).Value = 2
$ws.Cells.Item(26, 7).Value = 11.755292
$ws.Cells.Item(26, 8).Value = 23.510584
$ws.Cells.Item(26, 9).Value = 0.3793477806979606
$ws.Cells.Item(26, 10).Value = 0.2895063892140504
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 13).Value = 36.221842
$ws.Cells.Item(26, 14).Value = 72.443684
$ws.Cells.Item(26, 15).Value = 0.160199640292699
$ws.Cells.Item(26, 16).Value = 0.112824579196272
$ws.Cells.Item(26, 17).Value = 425.798329487864
$ws.Cells.Item(26, 18).Value = 1703.193317951456
$ws.Cells.Item(26, 19).Value = 0.06077137801364695
$ws.Cells.Item(26, 20).Value = 0.03266343653770736
